$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change), and a 4-row block
# reorder (OKB/InjectiveProtocol/Bittensor/dogwifhat) reflected by the diff.

$ws.Range("D2").Value = "'66.518.24"
$ws.Range("E2").Value = "  +4.36%  "
$ws.Range("D3").Value = "'3.488.35"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'588.36"
$ws.Range("E5").Value = "  +3.15%  "
$ws.Range("D6").Value = "'169.65"
$ws.Range("E6").Value = "  +7.46%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'3.487.68"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("D9").Value = "'0.579"
$ws.Range("E9").Value = "  +5.15%  "
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("E11").Value = "  +4.84%  "
$ws.Range("D12").Value = "'0.436"
$ws.Range("E12").Value = "  +3.53%  "
$ws.Range("D13").Value = "'4.093.64"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "'28.10"
$ws.Range("E15").Value = "  +4.10%  "
$ws.Range("D16").Value = "'66.543.92"
$ws.Range("E16").Value = "  +4.27%  "
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "'3.469.70"
$ws.Range("E18").Value = "  +3.15%  "
$ws.Range("E19").Value = "  +4.51%  "
$ws.Range("D20").Value = "'13.91"
$ws.Range("E20").Value = "  +2.68%  "
$ws.Range("D21").Value = "'386.83"
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("D22").Value = "'7.92"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("D23").Value = "'73.24"
$ws.Range("E23").Value = "  +2.85%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "'0.530"
$ws.Range("E25").Value = "  +3.33%  "
$ws.Range("D26").Value = "'0.0000121"
$ws.Range("E26").Value = "  +4.48%  "
$ws.Range("E27").Value = "  +5.84%  "
$ws.Range("E28").Value = "  +2.61%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'6.38"
$ws.Range("E30").Value = "  +7.14%  "
$ws.Range("D31").Value = "'1.49"
$ws.Range("E31").Value = "  +7.35%  "
$ws.Range("D32").Value = "'2.04"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("D33").Value = "'23.46"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("D34").Value = "'7.41"
$ws.Range("E34").Value = "  +6.54%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'1.60"
$ws.Range("E36").Value = "  +6.30%  "
$ws.Range("D37").Value = "'162.78"
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("D38").Value = "'0.875"
$ws.Range("E38").Value = "  +5.58%  "
$ws.Range("D39").Value = "'1.90"
$ws.Range("E39").Value = "  +4.68%  "
$ws.Range("D40").Value = "'0.0748"
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("D41").Value = "'4.64"
$ws.Range("E41").Value = "  +5.81%  "
$ws.Range("D42").Value = "'26.22"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("D43").Value = "'2.800.23"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "'6.59"
$ws.Range("E44").Value = "  +2.90%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'26.70"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'43.02"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.52"
$ws.Range("E47").Value = "  +7.01%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "'353.70"
$ws.Range("E48").Value = "  +6.16%  "
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("D50").Value = "'1.08"
$ws.Range("E50").Value = "  +3.72%  "
$ws.Range("D51").Value = "'33.70"
$ws.Range("E51").Value = "  +14.25%  "
